$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The FRIDAY column (F) held values that had gone through a pandas/numpy
# float pass ("2.0", "4.0\nIT-212-...", ..., "nan"). Fix them so they match
# the integer-style labels used by the rest of the schedule (same text as
# column C / D already use for these periods), and clear the three rows
# that only had a bogus "nan".

# Rows whose FRIDAY cell is plain text that does *not* contain a course
# line still need to end up as a shared string (not a numeric cell), so we
# go through Excel's "quote prefix" (leading apostrophe) input, then reset
# the cell style back to Normal so no visible/number formatting sticks.
$ws.Range("F2").Value = "'2"
$ws.Range("F2").Style = "Normal"

$ws.Range("F3").Value = "4`nIT-212-04009-Gilmartin-IDE-128A"
$ws.Range("F4").Value = "6`nIT-304-06008-Hogan-IDE-209"
$ws.Range("F5").Value = "8`nIT-140-08008-Fireheart-IDE-217A"

$ws.Range("F6").Value = "'10"
$ws.Range("F6").Style = "Normal"

$ws.Range("F7").Value = "12`nIT-200-12003-Pollitt-IDE-217A"

# Writing the multi-line labels above auto-expands the row height; put the
# affected rows back to their original auto-fit (default) height so only
# the cell contents differ, same as the source diff.
$ws.Rows("3:5").AutoFit()
$ws.Rows("7").AutoFit()

# The last three FRIDAY slots were "nan" placeholders - there's no class
# scheduled there, so just clear them out entirely.
$ws.Range("F8:F10").ClearContents()
